# Planning_TPI_Forestier.xlsx - "Mise a jour commentaire"
#
# 1. Sheet "Tâches réalisés": add the 09.mai journal rows (18-23)
# 2. Sheet "Planning effectif": mark a couple more half-day cells as done
#    (style change only)
# 3. View/selection bookkeeping on all three sheets + switch active tab
#    back to the first sheet

$wb = $excel.ActiveWorkbook

$wsPrev = $wb.Worksheets.Item(1)   # Planning prévisionel
$wsEff  = $wb.Worksheets.Item(2)   # Planning effectif
$wsTach = $wb.Worksheets.Item(3)   # Tâches réalisés

# ---------------------------------------------------------------------
# 1. Tâches réalisés - new "09.mai" block (rows 18-23)
# ---------------------------------------------------------------------

$wsTach.Cells.Item(18, 1).Value = 43594
$wsTach.Cells.Item(18, 1).NumberFormat = "d-mmm"
$wsTach.Cells.Item(18, 2).Value = "Test des différents moyen d'analyse"
$wsTach.Cells.Item(18, 3).Value = 0.0625
$wsTach.Cells.Item(18, 3).NumberFormat = "h:mm"
$wsTach.Cells.Item(18, 4).Value = "OpenCV" + [char]10 + "Scikit-image" + [char]10 + "Propre algorithme"
$wsTach.Rows.Item(18).RowHeight = 43.2

$wsTach.Cells.Item(19, 1).Value = 43594
$wsTach.Cells.Item(19, 1).NumberFormat = "d-mmm"
$wsTach.Cells.Item(19, 2).Value = "Visite des experts"
$wsTach.Cells.Item(19, 3).Value = 0.020833333333333332
$wsTach.Cells.Item(19, 3).NumberFormat = "h:mm"
$wsTach.Cells.Item(19, 4).Value = "Discuté avec M.Rulo sur l'avancée du projet et mon ressenti"

$wsTach.Cells.Item(20, 1).Value = 43594
$wsTach.Cells.Item(20, 1).NumberFormat = "d-mmm"
$wsTach.Cells.Item(20, 2).Value = "Implémentation de la méthode d'analyse avec scikit-image"
$wsTach.Cells.Item(20, 3).Value = 0.10416666666666667
$wsTach.Cells.Item(20, 3).NumberFormat = "h:mm"
$wsTach.Rows.Item(20).RowHeight = 28.8

$wsTach.Cells.Item(21, 1).Value = 43594
$wsTach.Cells.Item(21, 1).NumberFormat = "d-mmm"
$wsTach.Cells.Item(21, 2).Value = "Déplacement du robot en fonction de l'image analysé"
$wsTach.Cells.Item(21, 3).Value = 0.10416666666666667
$wsTach.Cells.Item(21, 3).NumberFormat = "h:mm"
$wsTach.Rows.Item(21).RowHeight = 28.8

$wsTach.Cells.Item(22, 1).Value = 43594
$wsTach.Cells.Item(22, 1).NumberFormat = "d-mmm"
$wsTach.Cells.Item(22, 2).Value = "Documentation"
$wsTach.Cells.Item(22, 3).Value = 0.041666666666666664
$wsTach.Cells.Item(22, 3).NumberFormat = "h:mm"

$wsTach.Cells.Item(23, 2).Value = "Total 09.mai"
$wsTach.Cells.Item(23, 3).Formula = "=SUM(C18:C22)"

# ---------------------------------------------------------------------
# 2. Planning effectif - highlight a couple more completed half-days
# ---------------------------------------------------------------------

$wsEff.Cells.Item(7, 7).Style  = $wsEff.Cells.Item(6, 7).Style   # G7  -> s24
$wsEff.Cells.Item(9, 8).Style  = $wsEff.Cells.Item(6, 7).Style   # H9  -> s24
$wsEff.Cells.Item(12, 7).Style = $wsEff.Cells.Item(6, 7).Style   # G12 -> s24
$wsEff.Cells.Item(12, 8).Style = $wsEff.Cells.Item(6, 7).Style   # H12 -> s24
$wsEff.Cells.Item(15, 7).Style = $wsEff.Cells.Item(6, 7).Style   # G15 -> s24
$wsEff.Cells.Item(15, 8).Style = $wsEff.Cells.Item(6, 7).Style   # H15 -> s24

# ---------------------------------------------------------------------
# 3. View bookkeeping - selection / scroll position on each sheet, and
#    make "Planning prévisionel" the active tab again
# ---------------------------------------------------------------------

$wsEff.Range("H15").Select()

$wsTach.Activate()
$wsTach.Range("C22").Select()

$wsPrev.Activate()
$wsPrev.Range("B3:X16").Select()
